$p = $ppt.ActivePresentation

# --- Slide 2 ("Agenda"): add three new agenda entries after "Problem Statement"
# and before "Project Approach" -------------------------------------------------
$agendaSlide = $p.Slides.Item(2)
$subtitle = $agendaSlide.Shapes.Item(2)
$tr = $subtitle.TextFrame.TextRange

# Paragraph 5 is "         Problem Statement" (same FF6600 / sz2800 run style
# used throughout this placeholder). Insert the three new lines right after it.
$problemStatementPara = $tr.Paragraphs(5, 1)
$problemStatementPara.InsertAfter("`r         nanoGPT Model Performance`r         nanoGPT Model Metrics`r         xLSTM Performance Metrics")

# --- Slide 7 ("nanoGPT Model"): extend the title to "nanoGPT Model Performance" ---
$nanoGptSlide = $p.Slides.Item(7)
$nanoGptTitle = $nanoGptSlide.Shapes.Item(1)
$titleRange = $nanoGptTitle.TextFrame.TextRange
$lastRun = $titleRange.Runs(3, 1)
$lastRun.Text = "l Performance"

# --- Slide 9 ("xLSTM Metrics"): rename title to "xLSTM Performance Metrics" ---
$xlstmSlide = $p.Slides.Item(9)
$xlstmTitle = $xlstmSlide.Shapes.Item(1)
$xlstmTitle.TextFrame.TextRange.Text = "xLSTM Performance Metrics"
